# The VIC "Chief Health Officer" press availability moved on: the "Last Date"
# (B6) is later, and the "News Link" (C6) now points at a YouTube recording of
# the presser instead of the old ABC News write-up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1) VIC row (row 6): updated "Last Date" -------------------------------
$ws.Range("B6").Value = 45087

# --- 2) Rebuild the "News Link" hyperlinks in column C ---------------------
# This engine's Hyperlink objects can only be swapped cleanly by clearing
# every hyperlink on the sheet and re-adding them (editing a single existing
# Hyperlink.Address in place leaves a stray duplicate relationship behind), so
# snapshot the current targets, swap in the one that changed (VIC), and
# recreate them all with the same displayed text they had before.

$newsLinks = @(
    @{ Cell = "C4";  Url = "https://www.health.gov.au/ministers/the-hon-mark-butler-mp/media/minister-for-health-and-aged-care-press-conference-30-march-2023" },
    @{ Cell = "C5";  Url = "https://www.youtube.com/watch?v=AOISAo2T3Rw" },
    @{ Cell = "C6";  Url = "https://www.youtube.com/watch?v=mvrt6kVMf1M" },
    @{ Cell = "C7";  Url = "https://www.abc.net.au/news/2023-04-04/qld-four-million-covid-19-coronavirus-chief-health-officer-2022/102180074" },
    @{ Cell = "C8";  Url = "https://www.news.com.au/national/western-australia/frustrated-wa-chief-health-officer-makes-two-big-calls-on-covid19/news-story/9f01a749a3562d3a2174896ec7dd48fa" },
    @{ Cell = "C9";  Url = "https://www.abc.net.au/news/2023-03-31/sa-covid-update/102174186" },
    @{ Cell = "C10"; Url = "https://www.premier.tas.gov.au/covid-19_updates/press_conference_-_06_july_2022" },
    @{ Cell = "C11"; Url = "https://www.abc.net.au/news/2023-02-17/act-covid-19-death-toll-far-higher-than-reported/101989422" },
    @{ Cell = "C12"; Url = "https://www.abc.net.au/news/2022-06-09/northern-territory-covid-natasha-fyles/101137616" }
)

$ws.Range("C4").Hyperlinks.Delete()

foreach ($link in $newsLinks) {
    $ws.Hyperlinks.Add($ws.Range($link.Cell), $link.Url, "", "", $link.Url) | Out-Null
    # Adding a hyperlink mints a fresh (duplicate) "Hyperlink" cell style;
    # re-apply the named style so the cell reuses the existing one.
    $ws.Range($link.Cell).Style = "Hyperlink"
}

# --- 3) Selection left where the workbook was last saved -------------------
$ws.Range("A4").Select()
